$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "261.13"
$ws.Range("E2").Value = "1.47%"
$rng.Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "5"
$ws.Range("G2").Style = "Normal"

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = "27.35"
$ws.Range("E3").Value = "1.20%"
$rng.Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "5"
$ws.Range("G3").Style = "Normal"

$rng = $ws.Range("D4:E4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = "4.758"
$ws.Range("E4").Value = "4.26%"
$rng.Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "5"
$ws.Range("G4").Style = "Normal"

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = "0.06065"
$ws.Range("E5").Value = "2.83%"
$rng.Style = "Normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "5"
$ws.Range("G5").Style = "Normal"

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "6.640"
$ws.Range("E6").Value = "0.16%"
$rng.Style = "Normal"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "5"
$ws.Range("G6").Style = "Normal"

$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = "0.8608"
$ws.Range("E7").Value = "1.05%"
$rng.Style = "Normal"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "5"
$ws.Range("G7").Style = "Normal"

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = "0.9214"
$ws.Range("E8").Value = "-1.87%"
$rng.Style = "Normal"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "5"
$ws.Range("G8").Style = "Normal"

$rng = $ws.Range("D9:E9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = "0.1410"
$ws.Range("E9").Value = "1.34%"
$rng.Style = "Normal"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "5"
$ws.Range("G9").Style = "Normal"

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = "0.04992"
$ws.Range("E10").Value = "-1.06%"
$rng.Style = "Normal"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "5"
$ws.Range("G10").Style = "Normal"

$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = "0.07094"
$ws.Range("E11").Value = "0.16%"
$rng.Style = "Normal"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "5"
$ws.Range("G11").Style = "Normal"

$rng = $ws.Range("D12:E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = "0.03045"
$ws.Range("E12").Value = "-0.78%"
$rng.Style = "Normal"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "5"
$ws.Range("G12").Style = "Normal"

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = "0.09089"
$ws.Range("E13").Value = "-0.31%"
$rng.Style = "Normal"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "5"
$ws.Range("G13").Style = "Normal"

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = "0.001541"
$ws.Range("E14").Value = "0.98%"
$rng.Style = "Normal"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "5"
$ws.Range("G14").Style = "Normal"

$rng = $ws.Range("D15:E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = "0.0006082"
$ws.Range("E15").Value = "0.35%"
$rng.Style = "Normal"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "5"
$ws.Range("G15").Style = "Normal"

$rng = $ws.Range("D16:E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = "0.006123"
$ws.Range("E16").Value = "0.03%"
$rng.Style = "Normal"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "5"
$ws.Range("G16").Style = "Normal"

$rng = $ws.Range("D17:E17")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = "3.452"
$ws.Range("E17").Value = "-1.11%"
$rng.Style = "Normal"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "5"
$ws.Range("G17").Style = "Normal"

$rng = $ws.Range("D18:E18")
$rng.NumberFormat = "@"
$ws.Range("D18").Value = "3.173"
$ws.Range("E18").Value = "-0.19%"
$rng.Style = "Normal"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "5"
$ws.Range("G18").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.27%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "5"
$ws.Range("G19").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.47%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "5"
$ws.Range("G20").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.23%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "5"
$ws.Range("G21").Style = "Normal"

$rng = $ws.Range("D22:E22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = "4.123"
$ws.Range("E22").Value = "5.03%"
$rng.Style = "Normal"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "5"
$ws.Range("G22").Style = "Normal"

$rng = $ws.Range("D23:E23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = "0.04246"
$ws.Range("E23").Value = "-0.49%"
$rng.Style = "Normal"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "5"
$ws.Range("G23").Style = "Normal"

$rng = $ws.Range("D24:E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = "0.001219"
$ws.Range("E24").Value = "-0.06%"
$rng.Style = "Normal"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "5"
$ws.Range("G24").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-8.75%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "5"
$ws.Range("G25").Style = "Normal"

$rng = $ws.Range("D26:E26")
$rng.NumberFormat = "@"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").Value = "0.04%"
$rng.Style = "Normal"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "5"
$ws.Range("G26").Style = "Normal"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.15%"
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "5"
$ws.Range("G27").Style = "Normal"

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "5"
$ws.Range("G28").Style = "Normal"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "5"
$ws.Range("G29").Style = "Normal"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "5"
$ws.Range("G30").Style = "Normal"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "5"
$ws.Range("G31").Style = "Normal"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "5"
$ws.Range("G32").Style = "Normal"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "5"
$ws.Range("G33").Style = "Normal"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "5"
$ws.Range("G34").Style = "Normal"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "5"
$ws.Range("G35").Style = "Normal"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "5"
$ws.Range("G36").Style = "Normal"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "5"
$ws.Range("G37").Style = "Normal"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "5"
$ws.Range("G38").Style = "Normal"

$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "5"
$ws.Range("G39").Style = "Normal"

$rng = $ws.Range("D40:E40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = "0.03879"
$ws.Range("E40").Value = "1.39%"
$rng.Style = "Normal"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "5"
$ws.Range("G40").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.13%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "5"
$ws.Range("G41").Style = "Normal"

$rng = $ws.Range("D42:E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = "0.004134"
$ws.Range("E42").Value = "-34.28%"
$rng.Style = "Normal"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "5"
$ws.Range("G42").Style = "Normal"

$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = "0.01500"
$ws.Range("E43").Value = "6.80%"
$rng.Style = "Normal"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "5"
$ws.Range("G43").Style = "Normal"

$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = "0.002161"
$ws.Range("E44").Value = "-11.44%"
$rng.Style = "Normal"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "5"
$ws.Range("G44").Style = "Normal"

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "0.00005318"
$ws.Range("E45").Value = "-0.47%"
$rng.Style = "Normal"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "5"
$ws.Range("G45").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.05%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "5"
$ws.Range("G46").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-18.41%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "5"
$ws.Range("G47").Style = "Normal"

$rng = $ws.Range("D48:E48")
$rng.NumberFormat = "@"
$ws.Range("D48").Value = "0.1321"
$ws.Range("E48").Value = "-47.51%"
$rng.Style = "Normal"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "5"
$ws.Range("G48").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.05%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "5"
$ws.Range("G49").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.05%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "5"
$ws.Range("G50").Style = "Normal"

$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "5"
$ws.Range("G51").Style = "Normal"
